{"js": "// Office.js (Word JavaScript API) script\n// Implements:\n//   1. Change the letter date \"September 19, 2025\" -> \"September 21, 2025\"\n//   2. Split the mailing-address paragraph \"989 Story Road, San Jose CA 95122\"\n//      into two paragraphs: \"989 Story Road\" and a new paragraph \"San Jose, CA 95122\"\n//   3. Remove the blank (NoSpacing) paragraph that immediately follows the\n//      \"...Board of Directors\" sign-off paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- 1. Update the date -------------------------------------------------\nconst dateResults = body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < dateResults.items.length; i++) {\n  dateResults.items[i].insertText(\"September 21, 2025\", Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// ---- 2. Split the mailing address paragraph -----------------------------\n// Locate the paragraph that holds the standalone mailing-address line\n// (the one right under the date, NOT the \"PROPERTY ADDRESS:\" table cell\n// later in the document which keeps its original combined text).\nconst addressCandidates = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"989 Story Road, San Jose CA 95122\") {\n    addressCandidates.push(paragraphs.items[i]);\n  }\n}\nfor (let i = 0; i < addressCandidates.length; i++) {\n  addressCandidates[i].parentTableOrNullObject.load(\"isNullObject\");\n}\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < addressCandidates.length; i++) {\n  if (addressCandidates[i].parentTableOrNullObject.isNullObject) {\n    addressParagraph = addressCandidates[i];\n    break;\n  }\n}\n\nif (addressParagraph) {\n  // Replacing the paragraph's own range with the two lines joined by a\n  // paragraph break splits it into two paragraphs in one step, both of\n  // which inherit the original paragraph/run formatting.\n  const addressRange = addressParagraph.getRange();\n  addressRange.insertText(\"989 Story Road\\nSan Jose, CA 95122\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// ---- 3. Remove the empty paragraph after \"Board of Directors\" -----------\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet boardParagraphIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(\"Board of Directors\") !== -1) {\n    boardParagraphIndex = i;\n    break;\n  }\n}\n\nif (boardParagraphIndex !== -1 && boardParagraphIndex + 1 < paragraphs.items.length) {\n  const nextParagraph = paragraphs.items[boardParagraphIndex + 1];\n  nextParagraph.load(\"text\");\n  await context.sync();\n  if (nextParagraph.text === \"\") {\n    nextParagraph.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Word COM interop (PowerShell-style) script\n# Implements:\n#   1. Change the letter date \"September 19, 2025\" -> \"September 21, 2025\"\n#   2. Split the mailing-address paragraph \"989 Story Road, San Jose CA 95122\"\n#      into two paragraphs: \"989 Story Road\" and a new paragraph \"San Jose, CA 95122\"\n#   3. Remove the blank (NoSpacing) paragraph that immediately follows the\n#      \"...Board of Directors\" sign-off paragraph.\n\n$d = $word.ActiveDocument\n\n# ---- 1. Update the date --------------------------------------------------\n# Assigning directly to the paragraph's Range.Text (rather than Find/Replace)\n# keeps the run's original xml:space=\"preserve\" formatting intact.\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq \"September 19, 2025\" + [char]13) {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# ---- 2. Split the mailing address paragraph ------------------------------\n# Locate the paragraph that holds the standalone mailing-address line (the\n# one right under the date). Skip the similarly-worded \"PROPERTY ADDRESS:\"\n# table cell further down, which keeps its original combined text.\n$addressParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -eq \"989 Story Road, San Jose CA 95122\" + [char]13 -and -not $p.Range.Information(12)) {\n        $addressParagraph = $p\n        break\n    }\n}\n\nif ($addressParagraph -ne $null) {\n    # Replacing the paragraph's text with the two lines joined by a paragraph\n    # mark (Chr 13) splits it into two paragraphs, inheriting the original\n    # paragraph/run formatting for both.\n    $addressParagraph.Range.Text = \"989 Story Road\" + [char]13 + \"San Jose, CA 95122\"\n}\n\n# ---- 3. Remove the empty paragraph after \"Board of Directors\" ------------\n$boardParagraph = $null\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs($i)\n    if ($p.Range.Text -like \"*Board of Directors*\") {\n        $boardParagraph = $p\n        break\n    }\n}\n\nif ($boardParagraph -ne $null) {\n    $nextParagraph = $boardParagraph.Next()\n    if ($nextParagraph -ne $null -and $nextParagraph.Range.Text -eq [char]13) {\n        $nextParagraph.Range.Delete()\n    }\n}\n"}
